$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: Md Zikrullah
$ws.Range("A2").Value = "18r21a12a0"
$ws.Range("B2").Value = "Md Zikrullah"
$ws.Range("C2").Value = "zikrullah@omnia.com"
$ws.Range("D2").Value = "ECE"
$ws.Range("E2").Value = "2022"

# Row 3: Vivek Jandhyala
$ws.Range("A3").Value = "18r21a1280"
$ws.Range("B3").Value = "Vivek Jandhyala"
$ws.Range("C3").Value = "vivek@omnia.com"
$ws.Range("D3").Value = "IT"
$ws.Range("E3").Value = "2022"

# Row 4: kota Srikar (previously the data in row 2)
$ws.Range("A4").Value = "18r21a1290"
$ws.Range("B4").Value = "kota Srikar"
$ws.Range("C4").Value = "kotasrikar009@gmail.com"
$ws.Range("D4").Value = "IT"
$ws.Range("E4").Value = "2022"

# Row 5: MD Q Arshad
$ws.Range("A5").Value = "18r21a1298"
$ws.Range("B5").Value = "MD Q Arshad"
$ws.Range("C5").Value = "arshad@omnia.com"
$ws.Range("D5").Value = "IT"
$ws.Range("E5").Value = "2022"
